$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Matching")

# Row 5 / new rows 6-7: replace "Function answer" with a red herring,
# and add a function definition plus a second red herring below it.
$ws.Range("C5").Value = "Red herring 1"
$ws.Range("C6").Value = "Function defintion"
$ws.Range("C7").Value = "Red herring 2"

# Row 4: change the "correct order" entry from D to E
$ws.Range("B4").Value = "E"

# Move the active selection down to reflect where the user left off editing.
$ws.Range("B9").Select()
